$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text so numeric-looking
# strings like "27.544.99" are not reinterpreted as numbers.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.544.99"
$ws.Range("D3").Value = "1.726.55"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "225.75"
$ws.Range("E5").Value = "  +3.33%  "
$ws.Range("D6").Value = "0.5364"
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.2670"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "0.06604"
$ws.Range("E9").Value = "  +4.17%  "
$ws.Range("D10").Value = "21.77"
$ws.Range("E10").Value = "  +6.77%  "
$ws.Range("D11").Value = "0.07712"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "4.612"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "1.725.91"
$ws.Range("E13").Value = "  +4.43%  "
$ws.Range("D14").Value = "1.963.01"
$ws.Range("E14").Value = "  +4.45%  "
$ws.Range("D15").Value = "0.5844"
$ws.Range("E15").Value = "  +4.42%  "
$ws.Range("D16").Value = "0.0₅8298"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "67.94"
$ws.Range("E17").Value = "  +3.87%  "
$ws.Range("D18").Value = "27.561.66"
$ws.Range("E18").Value = "  +5.56%  "
$ws.Range("D19").Value = "220.47"
$ws.Range("E19").Value = "  +15.37%  "
$ws.Range("D20").Value = "1.005"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "4.736"
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "6.096"
$ws.Range("E23").Value = "  +2.68%  "
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "148.35"
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("D26").Value = "1.709"
$ws.Range("E26").Value = "  +13.06%  "
$ws.Range("D27").Value = "0.1236"
$ws.Range("E27").Value = "  +3.84%  "
$ws.Range("D28").Value = "7.413"
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("D29").Value = "16.68"
$ws.Range("E29").Value = "  +4.63%  "
$ws.Range("D30").Value = "0.05569"
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("D31").Value = "1.303"
$ws.Range("E31").Value = "  +2.48%  "
$ws.Range("D32").Value = "3.553"
$ws.Range("E32").Value = "  +2.98%  "
$ws.Range("D33").Value = "3.461"
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("D34").Value = "1.661"
$ws.Range("E34").Value = "  +6.43%  "
$ws.Range("D35").Value = "0.9625"
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("D36").Value = "2.826"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").Value = "2.432"
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("D38").Value = "0.5952"
$ws.Range("E38").Value = "  +5.63%  "
$ws.Range("D39").Value = "0.01650"
$ws.Range("E39").Value = "  +4.65%  "
$ws.Range("D40").Value = "5.931"
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("D41").Value = "0.8566"
$ws.Range("E41").Value = "  +3.09%  "
$ws.Range("D42").Value = "1.057.80"
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "101.47"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "1.869.78"
$ws.Range("E45").Value = "  +4.39%  "
$ws.Range("D46").Value = "0.0₈114"
$ws.Range("E46").Value = "  +5.73%  "
$ws.Range("D48").Value = "8.217"
$ws.Range("E48").Value = "  +2.87%  "
$ws.Range("D49").Value = "0.4438"
$ws.Range("E49").Value = "  +2.31%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("B51").Value = "XinFinNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D51").Value = "0.06542"
$ws.Range("E51").Value = "  +12.87%  "

# Restore the original (default) cell style now that the text values are set,
# so no stray number-format style lingers on these cells.
$priceRange.Style = "Normal"
